$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New snapshot timestamp applied to rows 2-10 in column BH
$newTs = "2026-02-24 03:49:15"

# Row 2
$ws.Range("F2").Value = 3.35
$ws.Range("BH2").Value = $newTs

# Row 3
$ws.Range("BH3").Value = $newTs

# Row 4
$ws.Range("BH4").Value = $newTs

# Row 5
$ws.Range("F5").Value = 1.83
$ws.Range("G5").Value = 1.95
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.65
$ws.Range("Q5").Value = 2.32
$ws.Range("BH5").Value = $newTs

# Row 6
$ws.Range("BH6").Value = $newTs

# Row 7
$ws.Range("F7").Value = 2.1
$ws.Range("H7").Value = 3.65
$ws.Range("BH7").Value = $newTs

# Row 8
$ws.Range("G8").Value = 4.5
$ws.Range("H8").Value = 2.36
$ws.Range("P8").Value = 1.33
$ws.Range("Q8").Value = 3.5
$ws.Range("BH8").Value = $newTs

# Row 9
$ws.Range("F9").Value = 1.49
$ws.Range("G9").Value = 1.56
$ws.Range("H9").Value = 8.800000000000001
$ws.Range("I9").Value = 10.5
$ws.Range("K9").Value = 4.4
$ws.Range("P9").Value = 1.66
$ws.Range("Q9").Value = 2.26
$ws.Range("BH9").Value = $newTs

# Row 10
$ws.Range("BH10").Value = $newTs
